$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colB = @(0.142398651482722,0.1329491466732691,0.1272161570857122,0.1248974479733818,0.1245134918750779,0.1271848149585537,0.1391262300532361,0.1630850301430655,0.1810115073540857,0.1892358470044684,0.1923600379454484,0.191686754024289,0.1894926804650936,0.1881500194881056,0.1804754103396959,0.1757849566111247,0.1730936837499826,0.1721835958490772,0.1762835866173162,0.1901368682050872,0.1992478962401663,0.1943800036805925,0.1760581393712783,0.1565461964429744)
for ($i = 0; $i -lt $colB.Count; $i++) { $ws.Cells.Item($i + 2, 2).Value = $colB[$i] }

$colD = @(0.2510271085944709,0.2435338944786736,0.2390194271138881,0.2372016214916357,0.2369011019998339,0.2389948227732788,0.2484256048038702,0.2675993198005955,0.2820949850980128,0.2887770655859185,0.291319909301734,0.2907717096776992,0.2889860175654348,0.2878938507573992,0.2816600511684442,0.2778582245032197,0.275679800227195,0.274943651160612,0.2782620790409283,0.2895101810582048,0.2969341634379816,0.2929652450908975,0.2780794736656986,0.2623401255290787)
for ($i = 0; $i -lt $colD.Count; $i++) { $ws.Cells.Item($i + 2, 4).Value = $colD[$i] }

$colE = @(0.1815433724534969,0.1768847561558573,0.1741097564276686,0.1730004639147644,0.1728175697768606,0.1740947088093456,0.1799193922654041,0.1920171548577443,0.2013155274797072,0.2056343648307362,0.2072825385947041,0.2069270100871137,0.2057697064687005,0.2050624795570144,0.2010350655569795,0.1985871137473296,0.1971874941635576,0.1967150479634583,0.1988468358488191,0.2061092898233099,0.2109298552939052,0.2083502684457414,0.1987293913715718,0.1886722399632603)
for ($i = 0; $i -lt $colE.Count; $i++) { $ws.Cells.Item($i + 2, 5).Value = $colE[$i] }

$colF = @(1.019423977769264,1.024171066059395,1.027759155951287,1.029390591821397,1.029671711594119,1.027780472820766,1.020920979770231,1.012816052617943,1.010127270768429,1.009614717944359,1.009522904372453,1.009538127645996,1.009605114009823,1.00965946770711,1.010175074506847,1.010673450961043,1.0110269820577,1.011158165002861,1.010613475587938,1.009582661898705,1.009505170202971,1.009491944462567,1.010640381736643,1.014435588393589)
for ($i = 0; $i -lt $colF.Count; $i++) { $ws.Cells.Item($i + 2, 6).Value = $colF[$i] }

$colG = @(0.4819930937647854,0.4867705571386765,0.4901252131498453,0.4915980061757423,0.4918489434144249,0.4901446478526736,0.4835528137282381,0.4739776312218069,0.4689982852779551,0.4671820750247235,0.4665590887259441,0.4666903756994145,0.4671295221682783,0.4674069539058223,0.4691260310547776,0.470295764563275,0.4710108057412228,0.471260154699479,0.4701668710154152,0.4669987745437822,0.4653058881305014,0.4661747867035615,0.4702250112773569,0.4762078911293202)
for ($i = 0; $i -lt $colG.Count; $i++) { $ws.Cells.Item($i + 2, 7).Value = $colG[$i] }

$colH = @(0.6378099054947128,0.6442834492590208,0.6485948815336116,0.6504364781357879,0.6507473866652802,0.6486193751942793,0.639972111536288,0.6256860668090738,0.6168186502912434,0.6131382902108555,0.6117954704876638,0.6120824086710854,0.6130267963224156,0.613611884266291,0.6170662852865689,0.6192759981586278,0.6205802356056012,0.6210275412580586,0.6190373270551319,0.6127480264576235,0.6089340248056487,0.6109424989571437,0.6191451248395978,0.629264838606133)
for ($i = 0; $i -lt $colH.Count; $i++) { $ws.Cells.Item($i + 2, 8).Value = $colH[$i] }

$colJ = @(0.1770572214365203,0.1731240366063602,0.1708118169947141,0.1698954314000289,0.1697448289001926,0.1707993535542727,0.1756797393072063,0.1860654603550813,0.1941937512027039,0.1979998850095228,0.1994567720959424,0.1991423124927678,0.1981194317210253,0.1974949169711522,0.1939471936644992,0.1917985641451168,0.1905729461980172,0.1901597283940504,0.1920262321576445,0.1984194537395894,0.2026886358752193,0.200401787721546,0.1919232733657026,0.1831684726321328)
for ($i = 0; $i -lt $colJ.Count; $i++) { $ws.Cells.Item($i + 2, 10).Value = $colJ[$i] }

$colK = @(0.9074100319281229,0.7923363258557288,0.7213923435998879,0.692411714447843,0.6875953144851508,0.7210017828664661,0.8677936710947165,1.153286257948821,1.361510075855449,1.455888137311092,1.491575363373897,1.483891815330651,1.458825195695908,1.443464368897367,1.355335134382472,1.301181004883006,1.270000735564963,1.259438162028573,1.306949158005068,1.466189290264253,1.569959565145552,1.514603815202804,1.304341520151809,1.076315110838209)
for ($i = 0; $i -lt $colK.Count; $i++) { $ws.Cells.Item($i + 2, 11).Value = $colK[$i] }

$colO = @(2.200317181780179,2.223724351682023,2.239684896376389,2.246587967321375,2.247758304771253,2.239776378536433,2.208058046099083,2.158481067458126,2.129778912366234,2.118404173328571,2.114339198350109,2.11520387475872,2.118064884531151,2.119848917571716,2.130556169301698,2.137555897094003,2.141740243752096,2.143184161912259,2.136794378087103,2.117217954037557,2.105836713134664,2.111781627428286,2.137138162135713,2.170538818813213)
for ($i = 0; $i -lt $colO.Count; $i++) { $ws.Cells.Item($i + 2, 15).Value = $colO[$i] }
